$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price / 1h-volume figures (and the NEARProtocol <-> HuobiToken
# row swap at the bottom of the table) to match the latest GitHub Actions pull.

# Row 2
$ws.Range("D2").Value = '41.523.80'
$ws.Range("E2").Value = '  -0.66%  '

# Row 3
$ws.Range("D3").Value = '2.161.55'
$ws.Range("E3").Value = '  -2.52%  '

# Row 4
$ws.Range("E4").Value = '  +0.21%  '

# Row 5
$ws.Range("D5").Value = '237.53'
$ws.Range("E5").Value = '  -1.62%  '

# Row 6
$ws.Range("D6").Value = '0.605'
$ws.Range("E6").Value = '  -3.30%  '

# Row 7
$ws.Range("D7").Value = '70.89'
$ws.Range("E7").Value = '  -2.09%  '

# Row 8
$ws.Range("E8").Value = '  +0.08%  '

# Row 9
$ws.Range("D9").Value = '0.572'
$ws.Range("E9").Value = '  -3.79%  '

# Row 10
$ws.Range("D10").Value = '39.47'
$ws.Range("E10").Value = '  -6.15%  '

# Row 11
$ws.Range("D11").Value = '0.0899'
$ws.Range("E11").Value = '  -5.10%  '

# Row 12
$ws.Range("D12").Value = '54.06'
$ws.Range("E12").Value = '  -4.54%  '

# Row 13
$ws.Range("E13").Value = '  -3.65%  '

# Row 14
$ws.Range("E14").Value = '  -4.03%  '

# Row 15
$ws.Range("D15").Value = '2.487.73'
$ws.Range("E15").Value = '  -2.46%  '

# Row 16
$ws.Range("D16").Value = '14.19'
$ws.Range("E16").Value = '  -0.24%  '

# Row 17
$ws.Range("D17").Value = '2.172.13'
$ws.Range("E17").Value = '  -1.68%  '

# Row 18
$ws.Range("D18").Value = '0.781'
$ws.Range("E18").Value = '  -6.47%  '

# Row 19
$ws.Range("D19").Value = '41.443.71'
$ws.Range("E19").Value = '  -0.52%  '

# Row 20
$ws.Range("E20").Value = '  -4.39%  '

# Row 21
$ws.Range("D21").Value = '69.51'
$ws.Range("E21").Value = '  -4.01%  '

# Row 22
$ws.Range("E22").Value = '  -6.83%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.00'
$ws.Range("E23").Value = '  -10.20%  '

# Row 24
$ws.Range("E24").Value = '  -1.13%  '

# Row 25
$ws.Range("D25").Value = '1.97'
$ws.Range("E25").Value = '  -3.73%  '

# Row 26
$ws.Range("E26").Value = '  +0.03%  '

# Row 27
$ws.Range("D27").Value = '10.65'
$ws.Range("E27").Value = '  -6.56%  '

# Row 28
$ws.Range("D28").Value = '3.28'
$ws.Range("E28").Value = '  -9.53%  '

# Row 29
$ws.Range("D29").Value = '2.17'
$ws.Range("E29").Value = '  -4.75%  '

# Row 30
$ws.Range("E30").Value = '  -0.91%  '

# Row 31
$ws.Range("D31").Value = '171.53'
$ws.Range("E31").Value = '  +2.46%  '

# Row 32
$ws.Range("D32").Value = '19.69'
$ws.Range("E32").Value = '  -3.56%  '

# Row 33
$ws.Range("D33").Value = '32.81'
$ws.Range("E33").Value = '  +9.58%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0770'
$ws.Range("E34").Value = '  -3.38%  '

# Row 35
$ws.Range("D35").Value = '5.08'
$ws.Range("E35").Value = '  -7.88%  '

# Row 36
$ws.Range("E36").Value = '  -3.94%  '

# Row 37
$ws.Range("D37").Value = '0.104'
$ws.Range("E37").Value = '  -2.02%  '

# Row 38
$ws.Range("D38").Value = '4.23'
$ws.Range("E38").Value = '  -0.62%  '

# Row 39
$ws.Range("D39").Value = '0.0298'
$ws.Range("E39").Value = '  -1.29%  '

# Row 40
$ws.Range("D40").Value = '12.02'
$ws.Range("E40").Value = '  -9.53%  '

# Row 41
$ws.Range("E41").Value = '  -2.79%  '

# Row 42
$ws.Range("D42").Value = '5.32'
$ws.Range("E42").Value = '  -5.26%  '

# Row 43
$ws.Range("D43").Value = '58.36'
$ws.Range("E43").Value = '  -8.64%  '

# Row 44
$ws.Range("D44").Value = '0.188'
$ws.Range("E44").Value = '  -4.53%  '

# Row 45
$ws.Range("D45").Value = '8.33'
$ws.Range("E45").Value = '  -4.34%  '

# Row 46
$ws.Range("D46").Value = '0.0958'
$ws.Range("E46").Value = '  -3.99%  '

# Row 47
$ws.Range("D47").Value = '95.39'
$ws.Range("E47").Value = '  -7.25%  '

# Row 48
$ws.Range("D48").Value = '1.07'
$ws.Range("E48").Value = '  -2.94%  '

# Row 49
$ws.Range("E49").Value = '  -4.84%  '

# Row 50
$ws.Range("B50").Value = 'HuobiToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D50").Value = '2.62'
$ws.Range("E50").Value = '  -2.47%  '

# Row 51
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '2.15'
$ws.Range("E51").Value = '  -7.89%  '
